$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bondora: switch to English website -- translate header row from German to English
$ws.Range("A1").Value = "Period"
$ws.Range("B1").Value = "Opening balance"
$ws.Range("C1").Value = "Net capital deployed"
$ws.Range("D1").Value = "Net loan investments"
$ws.Range("E1").Value = "Principal received - total"
$ws.Range("F1").Value = "Interest received - total"
$ws.Range("G1").Value = "Principal and interest received - total"
$ws.Range("H1").Value = "Closing balance"
$ws.Range("I1").Value = "Principal planned - total"
$ws.Range("J1").Value = "Interest planned - total"
$ws.Range("K1").Value = "Principal and interest planned - total"

# Header row no longer has the bottom border, and is now centred with wrapped text
$headerRow = $ws.Range("A1:K1")
$headerRow.Borders.LineStyle = -4142
$headerRow.HorizontalAlignment = -4108
$headerRow.WrapText = $true

# Header row height shrinks slightly to match the other rows
$ws.Rows.Item(1).RowHeight = 13.8

# Move/restore the active selection
$r1 = $ws.Range("E1:G1")
$r2 = $ws.Range("I1:K1")
$r3 = $ws.Range("D6")
$u = $excel.Union($r1, $r2, $r3)
$u.Select()
